# "optimized code and applied ngxSpinnerservies"
# Rename Sheet2 -> Modules and populate it with the Modules master list.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Modules"

# --- Data: row number | module name | created-date (OLE serial) ---------
$moduleRows = @(
    @(1,  "Training of Trainer & Training of Assessor",                                 44043.942811365741, 90),
    @(2,  "Qualification Pack – National Occupational Standard Builder",                44043.943102164354, 120),
    @(3,  "Center Accreditation and Affiliation Module",                                 44043.943250995369, 105),
    @(4,  "Candidate Self Registration and Candidate Login and Profile Management",      44043.944480011574, 150),
    @(5,  "Fee Based Module",                                                            44043.944788506946, 45),
    @(6,  "NON PMKVY Module",                                                            44043.944943425930, 45),
    @(7,  "Rozgar Mela",                                                                 44043.945096261574, 30),
    @(8,  "TP – TC Onboarding and Target Allocation",                                    44043.945341944447, 90),
    @(9,  "Candidate Training Lifecycle",                                                44043.945489398146, 60),
    @(10, "Batch Creation and Candidate Enrollment",                                     44043.945682604164, 105),
    @(11, "Assessment, Re-Assessment and Certification",                                 44043.945906215275, 90),
    @(12, "Budget and Payout",                                                           44043.946069409722, 45),
    @(13, "Placement",                                                                   44043.946205601853, 30),
    @(14, "Continuous Monitoring",                                                       44043.946392812497, 60),
    @(15, "Third Party Integrations",                                                    44043.946536469906, 60),
    @(16, "API Integration - States",                                                    44043.946729618059, 60),
    @(17, "API Integration - Central Ministry",                                          44043.946926041666, 75),
    @(18, "Nano BI Report",                                                              44043.947062581021, 30),
    @(19, "Nano BI Dashboard",                                                           44043.947199895832, 45),
    @(20, "Excel Report",                                                                44043.947426550927, 30)
)

$lastRow = $moduleRows.Count

# --- Formatting shared by the whole table: thin black border, wrapped, --
# --- vertically centred text (matches sheet "User"'s existing style). ---
$tableRange = $ws2.Range("A1:D$lastRow")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.ColorIndex = 1
$tableRange.WrapText = $true
$tableRange.VerticalAlignment = -4108

foreach ($row in $moduleRows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $r
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 3).NumberFormat = "mm:ss.0"
    $ws2.Cells.Item($r, 4).Value = 1
    $ws2.Rows.Item($r).RowHeight = $row[3]
}

# Column C (created-date) was narrowed slightly from the sheet default.
$ws2.Columns.Item(3).ColumnWidth = 8.3

# Page setup picked up from the print dialog (A4-ish letter, portrait).
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection left on F1 after data entry.
[void]$ws2.Range("F1").Select()
